$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E hold numeric-looking values (prices / percentages) that
# must remain stored as text, matching the original workbook formatting.
# Force each cell to Text format before writing so Excel does not silently
# convert the string into a numeric value.
$numericLookingCells = @(
    'D2', 'E2', 'E3', 'D4', 'E4', 'E5', 'D6', 'E6', 'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'D10', 'E10', 'D11', 'E11', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'D18', 'E18', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'D28', 'E28', 'D40', 'E40', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'D45', 'E45', 'E46', 'D47', 'E47', 'D48', 'E48', 'E49', 'E50'
)
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Updated price / volume values
$ws.Range('D2').Value = '259.59'
$ws.Range('E2').Value = '5.85%'
$ws.Range('E3').Value = '-3.47%'
$ws.Range('D4').Value = '5.222'
$ws.Range('E4').Value = '-0.74%'
$ws.Range('E5').Value = '3.87%'
$ws.Range('D6').Value = '6.728'
$ws.Range('E6').Value = '1.72%'
$ws.Range('D7').Value = '0.8737'
$ws.Range('E7').Value = '2.58%'
$ws.Range('D8').Value = '0.9941'
$ws.Range('E8').Value = '16.12%'
$ws.Range('D9').Value = '0.1425'
$ws.Range('E9').Value = '3.93%'
$ws.Range('D10').Value = '0.07232'
$ws.Range('E10').Value = '2.60%'
$ws.Range('D11').Value = '0.03202'
$ws.Range('E11').Value = '0.33%'
$ws.Range('D12').Value = '0.09246'
$ws.Range('E12').Value = '-0.02%'
$ws.Range('D13').Value = '0.001545'
$ws.Range('E13').Value = '1.22%'
$ws.Range('D14').Value = '0.0006071'
$ws.Range('E14').Value = '-93.92%'
$ws.Range('D15').Value = '0.005947'
$ws.Range('E15').Value = '-0.82%'
$ws.Range('D16').Value = '3.498'
$ws.Range('E16').Value = '0.16%'
$ws.Range('D17').Value = '3.233'
$ws.Range('E17').Value = '1.28%'
$ws.Range('D18').Value = '2.207'
$ws.Range('E18').Value = '1.48%'
$ws.Range('E19').Value = '-1.15%'
$ws.Range('D20').Value = '0.03637'
$ws.Range('E20').Value = '10.97%'
$ws.Range('D21').Value = '0.1290'
$ws.Range('E21').Value = '1.01%'
$ws.Range('D22').Value = '3.516'
$ws.Range('E22').Value = '0.72%'
$ws.Range('D23').Value = '0.04167'
$ws.Range('E23').Value = '1.87%'
$ws.Range('D24').Value = '0.1397'
$ws.Range('E24').Value = '1.24%'
$ws.Range('D25').Value = '0.001217'
$ws.Range('E25').Value = '-0.52%'
$ws.Range('D26').Value = '0.004574'
$ws.Range('E26').Value = '10.50%'
$ws.Range('D27').Value = '0.0001197'
$ws.Range('E27').Value = '-0.29%'
$ws.Range('D28').Value = '0.0001935'
$ws.Range('E28').Value = '33.49%'
$ws.Range('D40').Value = '0.03847'
$ws.Range('E40').Value = '2.42%'
$ws.Range('D41').Value = '0.005487'
$ws.Range('E41').Value = '6.78%'
$ws.Range('D42').Value = '0.1109'
$ws.Range('E42').Value = '4.32%'
$ws.Range('D43').Value = '0.002371'
$ws.Range('E43').Value = '-1.31%'
$ws.Range('D44').Value = '0.01091'
$ws.Range('E44').Value = '16.48%'
$ws.Range('D45').Value = '0.00005414'
$ws.Range('E45').Value = '2.25%'
$ws.Range('E46').Value = '-0.29%'
$ws.Range('D47').Value = '0.08538'
$ws.Range('E47').Value = '13.77%'
$ws.Range('D48').Value = '0.002136'
$ws.Range('E48').Value = '-12.53%'
$ws.Range('E49').Value = '-0.29%'
$ws.Range('E50').Value = '-0.29%'

# Coin name / link swap between rows 41 and 42, plus link text updates
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
